# EvaChecks.xlsx - "added 2 new issues to the evachecks"
#
# Inserts two new High-Priority issue rows (Overlapping CE Participation
# Records / Overlapping HMIS Participation Records) just above the existing
# "Missing Address" row, renames the sheet to "EvaChecks (11)", and keeps
# the AutoFilter / _FilterDatabase defined name / selection in sync with
# the new, larger data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new blank rows right before the old row 110 (pushes the
#    existing "Error"/"Warning" rows down by two, from 110-117 to 112-119).
# ---------------------------------------------------------------------------
$ws.Rows("110:111").Insert()

# Fully clear any formatting inherited from the row above so that re-entering
# values below results in the plain/default cell style used elsewhere in the
# sheet for rows that don't use a custom row height.
$ws.Range("A110:G111").Clear()

# ---------------------------------------------------------------------------
# 2. Populate the first new row (110) - Overlapping CE Participation Records.
# ---------------------------------------------------------------------------
$ws.Range("A110").Value = "pdde"
$ws.Range("B110").Value = "High Priority"
$ws.Range("C110").Value = "Overlapping CE Participation Records"
$ws.Range("E110").Value = "This project has more than one CE Participation record that covers the same time period. Please be sure you are ending any records that are no longer accurate before creating a new record."
$ws.Range("D110").Value = "2.09 CE Participation Status"
$ws.Range("G110").Value = 128
$ws.Range("G110").NumberFormat = "#,##0"
$ws.Range("G110").HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# 3. Populate the second new row (111) - Overlapping HMIS Participation Records.
# ---------------------------------------------------------------------------
$ws.Range("A111").Value = "pdde"
$ws.Range("B111").Value = "High Priority"
$ws.Range("C111").Value = "Overlapping HMIS Participation Records"
$ws.Range("D111").Value = "2.08 HMIS Participation"
$ws.Range("E111").Value = "This project has more than one HMIS Participation record that covers the same time period. Please be sure you are ending any records that are no longer accurate before creating a new record."
$ws.Range("G111").Value = 129
$ws.Range("G111").NumberFormat = "#,##0"
$ws.Range("G111").HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# 4. Refresh the AutoFilter range so it covers the two new rows
#    (A1:G119 instead of A1:G117).
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:G119").AutoFilter()

# ---------------------------------------------------------------------------
# 5. Rename the worksheet and keep the hidden _FilterDatabase defined name
#    pointing at the (renamed) sheet and the new, larger range.
# ---------------------------------------------------------------------------
$ws.Name = "EvaChecks (11)"
$wb.Names.Item(1).RefersTo = "='EvaChecks (11)'!`$A`$1:`$G`$119"

# ---------------------------------------------------------------------------
# 6. Leave the active selection on the first freshly-added row, matching
#    where the author was last working in the sheet.
# ---------------------------------------------------------------------------
$ws.Range("A112").Select()
